# Insert a new model run "ibes_2|fwdepsqcut|tune_mse_ind3" as a new row 13
# in the "average_mae" sheet, pushing existing rows 13..60 down to 14..61.
# (mirrors: select row 13, Insert row, fill in the new metrics)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 13..60 down to 14..61 (process bottom-up so nothing is
# clobbered before it has been copied to its new location).
for ($r = 60; $r -ge 13; $r--) {
    $src = $ws.Range("A" + $r + ":I" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":I" + ($r + 1))
    $dst.Value2 = $src.Value2
}

# Column A carries the bold/centered/bordered "label" style. Row 61 is a
# brand-new cell (previously past the end of the sheet) and row 13 is about
# to get brand-new content, so both need that formatting applied; copy it
# from an existing label cell (A12) rather than re-building it by hand.
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A61").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null

# Fill the new row 13 with the new model run's results.
$ws.Range("A13").Value2 = "ibes_2|fwdepsqcut|tune_mse_ind3"
$ws.Range("B13").Value2 = 0.00928508302271853
$ws.Range("C13").Value2 = 0.009674189396799985
$ws.Range("D13").Value2 = 0.0003512206687837434
$ws.Range("E13").Value2 = 0.0003384092866814527
$ws.Range("F13").Value2 = 0.2301836238080931
$ws.Range("G13").Value2 = 0.2582640092197721
$ws.Range("H13").Value2 = 0.2582640092197721
$ws.Range("I13").Value2 = 14156

Write-Host "inserted new row 13 (tune_mse_ind3), shifted rows 13-60 to 14-61"
